$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "Primera" (row2/row4) and "Segunda" (row3/row5) quality
# got swapped between the two reporting dates (44195 and 44216).
# Swap row 2 <-> row 4, and row 3 <-> row 5 for columns D, N, O, P, S.

function Swap-Row($r1, $r2) {
    $cols = @("D", "N", "O", "P", "S")
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $tmp = $cell1.Value()
        $cell1.Value = $cell2.Value()
        $cell2.Value = $tmp
    }
}

Swap-Row 2 4
Swap-Row 3 5
